# Update PollsData workbook with new elabe / harris poll rows (167-169)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Row 167 - elabe poll (id 66)
# ---------------------------------------------------------------
$ws.Cells.Item(167, 1).Value  = 66        # A id
$ws.Cells.Item(167, 2).Value  = 2022      # B year
$ws.Cells.Item(167, 3).Value  = 20        # C week
$ws.Cells.Item(167, 4).Value  = 1         # D month
$ws.Cells.Item(167, 5).Value  = 11        # E day
$ws.Cells.Item(167, 6).Value  = "elabe"      # F firm
$ws.Cells.Item(167, 7).Value  = "online"     # G collectmode
$ws.Cells.Item(167, 8).Value  = "partially"  # H unsure
$ws.Cells.Item(167, 9).Value  = 0         # I rolling
$ws.Cells.Item(167, 10).Value = 943       # J n
$ws.Cells.Item(167, 11).Value = 1.5       # K c_poutou
$ws.Cells.Item(167, 12).Value = 1         # L c_arthaud
$ws.Cells.Item(167, 13).Value = 9.5       # M c_melenchon
$ws.Cells.Item(167, 14).Value = 1.5       # N c_roussel
$ws.Cells.Item(167, 15).Value = 1         # O c_montebourg
$ws.Cells.Item(167, 16).Value = 7         # P c_jadot
$ws.Cells.Item(167, 17).Value = 3.5       # Q c_hidalgo
$ws.Cells.Item(167, 18).Value = 23        # R c_macron
$ws.Cells.Item(167, 19).Value = 17        # S c_pecresse
$ws.Cells.Item(167, 22).Value = 2         # V c_lassalle
$ws.Cells.Item(167, 23).Value = 1.5       # W c_daignant
$ws.Cells.Item(167, 24).Value = 17        # X c_lepen
$ws.Cells.Item(167, 25).Value = 13        # Y c_zemmour
$ws.Cells.Item(167, 26).Value = "T_0.5"      # Z c_asselineau
$ws.Cells.Item(167, 28).Value = 1         # AB c_philippot
$ws.Cells.Item(167, 31).Value = 0.5       # AE c_thouy
$ws.Cells.Item(167, 32).Value = 1         # AF omit

# ---------------------------------------------------------------
# Row 168 - elabe poll (id 66)
# ---------------------------------------------------------------
$ws.Cells.Item(168, 1).Value  = 66        # A id
$ws.Cells.Item(168, 2).Value  = 2022      # B year
$ws.Cells.Item(168, 3).Value  = 20        # C week
$ws.Cells.Item(168, 4).Value  = 1         # D month
$ws.Cells.Item(168, 5).Value  = 11        # E day
$ws.Cells.Item(168, 6).Value  = "elabe"      # F firm
$ws.Cells.Item(168, 7).Value  = "online"     # G collectmode
$ws.Cells.Item(168, 8).Value  = "partially"  # H unsure
$ws.Cells.Item(168, 9).Value  = 0         # I rolling
$ws.Cells.Item(168, 10).Value = 953       # J n
$ws.Cells.Item(168, 11).Value = 1         # K c_poutou
$ws.Cells.Item(168, 12).Value = 0.5       # L c_arthaud
$ws.Cells.Item(168, 13).Value = 10        # M c_melenchon
$ws.Cells.Item(168, 14).Value = 1.5       # N c_roussel
$ws.Cells.Item(168, 15).Value = 1         # O c_montebourg
$ws.Cells.Item(168, 16).Value = 6         # P c_jadot
$ws.Cells.Item(168, 17).Value = 2.5       # Q c_hidalgo
$ws.Cells.Item(168, 18).Value = 22.5      # R c_macron
$ws.Cells.Item(168, 19).Value = 16        # S c_pecresse
$ws.Cells.Item(168, 22).Value = 2         # V c_lassalle
$ws.Cells.Item(168, 23).Value = 2         # W c_daignant
$ws.Cells.Item(168, 24).Value = 16.5      # X c_lepen
$ws.Cells.Item(168, 25).Value = 12.5      # Y c_zemmour
$ws.Cells.Item(168, 26).Value = 0.5       # Z c_asselineau
$ws.Cells.Item(168, 28).Value = 1         # AB c_philippot
$ws.Cells.Item(168, 30).Value = 4         # AD c_taubira
$ws.Cells.Item(168, 31).Value = 0.5       # AE c_thouy

# ---------------------------------------------------------------
# Row 169 - harris poll (id 67)
# ---------------------------------------------------------------
$ws.Cells.Item(169, 1).Value  = 67        # A id
$ws.Cells.Item(169, 2).Value  = 2022      # B year
$ws.Cells.Item(169, 3).Value  = 19        # C week
$ws.Cells.Item(169, 4).Value  = 1         # D month
$ws.Cells.Item(169, 5).Value  = 9         # E day
$ws.Cells.Item(169, 6).Value  = "harris"     # F firm
$ws.Cells.Item(169, 7).Value  = "online"     # G collectmode
$ws.Cells.Item(169, 8).Value  = "included"   # H unsure
$ws.Cells.Item(169, 9).Value  = 0         # I rolling
$ws.Cells.Item(169, 10).Value = 2124      # J n
$ws.Cells.Item(169, 11).Value = 1         # K c_poutou
$ws.Cells.Item(169, 12).Value = "T_0.5"      # L c_arthaud
$ws.Cells.Item(169, 13).Value = 11        # M c_melenchon
$ws.Cells.Item(169, 14).Value = 2         # N c_roussel
$ws.Cells.Item(169, 15).Value = 1         # O c_montebourg
$ws.Cells.Item(169, 16).Value = 7         # P c_jadot
$ws.Cells.Item(169, 17).Value = 4         # Q c_hidalgo
$ws.Cells.Item(169, 18).Value = 25        # R c_macron
$ws.Cells.Item(169, 19).Value = 16        # S c_pecresse
$ws.Cells.Item(169, 22).Value = "T_0.5"      # V c_lassalle
$ws.Cells.Item(169, 23).Value = 2         # W c_daignant
$ws.Cells.Item(169, 24).Value = 16        # X c_lepen
$ws.Cells.Item(169, 25).Value = 15        # Y c_zemmour
$ws.Cells.Item(169, 26).Value = "T_0.5"      # Z c_asselineau

$ab169 = $ws.Cells.Item(169, 28)          # AB c_philippot (styled, black font)
$ab169.Value = "T_0.5"
$ab169.Font.Color = 0

# ---------------------------------------------------------------
# Update selection / active cell to mirror the authored view state
# ---------------------------------------------------------------
$ws.Range("AC169").Select()
